$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (idx=2, "닌자 백팩"): merit text now grants two shuriken, and
# the demerit (MoveSpeed penalty) is removed.
$ws.Range("C4").Value = "수리검을 두 개 더 가질 수 있습니다."
$ws.Range("E4").Value = "none"
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = "none"
$ws.Range("H4").Value = 0

# Row 7 (idx=5): "차징속도 감소" (charge speed decrease) becomes
# "차징 속도 증가" (charge speed increase) -- now a merit instead of a demerit.
$ws.Range("B7").Value = "차징 속도 증가"
$ws.Range("C7").Value = "수리검의 차징 속도가 100% 증가합니다."
$ws.Range("E7").Value = "ChargeSpeed"
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = "none"
$ws.Range("H7").Value = 0

# Row 8 (idx=6, "바운스!"): bounced shuriken damage multiplier tweak.
$ws.Range("C8").Value = "수리검이 벽에 부딪히면 튕깁니다. 튕긴 수리검은 1.5배의 데미지를 줍니다.(중첩불가)"

# Row 13 (idx=11, "거미줄"): spiderweb effect description updated.
$ws.Range("C13").Value = "수리검이 거미줄을 소환합니다. 거미줄 안에서는 느려지고 대시할 수 없습니다.(중첩불가)"

# Restore the active selection cursor to C13, as in the saved workbook.
$ws.Range("C13").Select()
